# Autogenerated on Sun Feb 01 2015 22:24:41 GMT-0500 (Eastern Standard Time)
#
# - Rename the single worksheet "Data" -> "Summary"
# - Remove the breakdown row (Micro / SMEs / MSMEs) that lived in row 5,
#   shrinking the used range back down to A1:A3
# - Register a new named cell style "title_" (bold + underlined Calibri 11),
#   inserted alongside the workbook's existing "name"/"title"/"source"/
#   "HyperLink" styles (not applied to any cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Summary"

$ws.Rows("5:5").Delete()

$titleStyle = $wb.Styles.Add("title_")
$titleStyle.Font.Bold = $true
$titleStyle.Font.Underline = $true
